$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recomputed TPM-based NATMI values (Cd34-Sele LR pairs) per new TPM inputs.
$ws.Range("G2").Value = 171.4863125
$ws.Range("H2").Value = 342.972625
$ws.Range("I2").Value = 0.6883220549653215
$ws.Range("J2").Value = 0.6156469349130842
$ws.Range("M2").Value = 7.369448
$ws.Range("N2").Value = 14.738896
$ws.Range("O2").Value = 0.7452608427984224
$ws.Range("P2").Value = 0.661061693471796
$ws.Range("Q2").Value = 1263.7594626805
$ws.Range("R2").Value = 5055.037850722
$ws.Range("S2").Value = 0.5129794748001976
$ws.Range("T2").Value = 0.406980605374364
$ws.Range("G3").Value = 171.4863125
$ws.Range("H3").Value = 342.972625
$ws.Range("I3").Value = 0.6883220549653215
$ws.Range("J3").Value = 0.6156469349130842
$ws.Range("O3").Value = 0.01116592909756377
$ws.Range("P3").Value = 0.01485661309677453
$ws.Range("Q3").Value = 18.93437538416666
$ws.Range("R3").Value = 113.606252305
$ws.Range("S3").Value = 0.007685755262032173
$ws.Range("T3").Value = 0.009146428316218822
$ws.Range("G4").Value = 171.4863125
$ws.Range("H4").Value = 342.972625
$ws.Range("I4").Value = 0.6883220549653215
$ws.Range("J4").Value = 0.6156469349130842
$ws.Range("M4").Value = 0.084843
$ws.Range("N4").Value = 0.254529
$ws.Range("O4").Value = 0.008580040959044227
$ws.Range("P4").Value = 0.0114160091622658
$ws.Range("Q4").Value = 14.5494132114375
$ws.Range("R4").Value = 87.29647926862501
$ws.Range("S4").Value = 0.00590583142461595
$ws.Range("T4").Value = 0.007028231049688627
$ws.Range("G5").Value = 171.4863125
$ws.Range("H5").Value = 342.972625
$ws.Range("I5").Value = 0.6883220549653215
$ws.Range("J5").Value = 0.6156469349130842
$ws.Range("M5").Value = 2.32371
$ws.Range("N5").Value = 6.97113
$ws.Range("O5").Value = 0.2349931871449696
$ws.Range("P5").Value = 0.3126656842691638
$ws.Range("Q5").Value = 398.484459219375
$ws.Range("R5").Value = 2390.90675531625
$ws.Range("S5").Value = 0.1617509934784759
$ws.Range("T5").Value = 0.1924916701728128
$ws.Range("I6").Value = 0.2310851226340305
$ws.Range("J6").Value = 0.310029686890149
$ws.Range("M6").Value = 7.369448
$ws.Range("N6").Value = 14.738896
$ws.Range("O6").Value = 0.7452608427984224
$ws.Range("P6").Value = 0.661061693471796
$ws.Range("Q6").Value = 424.272342149712
$ws.Range("R6").Value = 2545.634052898272
$ws.Range("S6").Value = 0.1722186932524143
$ws.Range("T6").Value = 0.2049487498421325
$ws.Range("I7").Value = 0.2310851226340305
$ws.Range("J7").Value = 0.310029686890149
$ws.Range("O7").Value = 0.01116592909756377
$ws.Range("P7").Value = 0.01485661309677453
$ws.Range("S7").Value = 0.002580280094833413
$ws.Range("T7").Value = 0.004605991106641093
$ws.Range("I8").Value = 0.2310851226340305
$ws.Range("J8").Value = 0.310029686890149
$ws.Range("M8").Value = 0.084843
$ws.Range("N8").Value = 0.254529
$ws.Range("O8").Value = 0.008580040959044227
$ws.Range("P8").Value = 0.0114160091622658
$ws.Range("Q8").Value = 4.884563718342
$ws.Range("R8").Value = 43.961073465078
$ws.Range("S8").Value = 0.001982719817225739
$ws.Range("T8").Value = 0.003539301746112338
$ws.Range("I9").Value = 0.2310851226340305
$ws.Range("J9").Value = 0.310029686890149
$ws.Range("M9").Value = 2.32371
$ws.Range("N9").Value = 6.97113
$ws.Range("O9").Value = 0.2349931871449696
$ws.Range("P9").Value = 0.3126656842691638
$ws.Range("Q9").Value = 133.78015343574
$ws.Range("R9").Value = 1204.02138092166
$ws.Range("S9").Value = 0.05430342946955698
$ws.Range("T9").Value = 0.09693564419526303
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.04101466666666666
$ws.Range("H10").Value = 0.123044
$ws.Range("I10").Value = 0.0001646271310645721
$ws.Range("J10").Value = 0.0002208679525354116
$ws.Range("M10").Value = 7.369448
$ws.Range("N10").Value = 14.738896
$ws.Range("O10").Value = 0.7452608427984224
$ws.Range("P10").Value = 0.661061693471796
$ws.Range("Q10").Value = 0.3022554532373333
$ws.Range("R10").Value = 1.813532719424
$ws.Range("S10").Value = 0.0001226901544446694
$ws.Range("T10").Value = 0.0001460073427367075
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.04101466666666666
$ws.Range("H11").Value = 0.123044
$ws.Range("I11").Value = 0.0001646271310645721
$ws.Range("J11").Value = 0.0002208679525354116
$ws.Range("O11").Value = 0.01116592909756377
$ws.Range("P11").Value = 0.01485661309677453
$ws.Range("Q11").Value = 0.004528566062222222
$ws.Range("R11").Value = 0.04075709456
$ws.Range("S11").Value = 0.000001838214873002351
$ws.Range("T11").Value = 0.000003281349716295371
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.04101466666666666
$ws.Range("H12").Value = 0.123044
$ws.Range("I12").Value = 0.0001646271310645721
$ws.Range("J12").Value = 0.0002208679525354116
$ws.Range("M12").Value = 0.084843
$ws.Range("N12").Value = 0.254529
$ws.Range("O12").Value = 0.008580040959044227
$ws.Range("P12").Value = 0.0114160091622658
$ws.Range("Q12").Value = 0.003479807364
$ws.Range("R12").Value = 0.031318266276
$ws.Range("S12").Value = 0.000001412507527503971
$ws.Range("T12").Value = 0.000002521430569795147
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.04101466666666666
$ws.Range("H13").Value = 0.123044
$ws.Range("I13").Value = 0.0001646271310645721
$ws.Range("J13").Value = 0.0002208679525354116
$ws.Range("M13").Value = 2.32371
$ws.Range("N13").Value = 6.97113
$ws.Range("O13").Value = 0.2349931871449696
$ws.Range("P13").Value = 0.3126656842691638
$ws.Range("Q13").Value = 0.09530619108000001
$ws.Range("R13").Value = 0.8577557197200001
$ws.Range("S13").Value = 0.00003868625421939644
$ws.Range("T13").Value = 0.00006905782951261367
$ws.Range("G14").Value = 18.830862
$ws.Range("H14").Value = 37.661724
$ws.Range("I14").Value = 0.07558444426057843
$ws.Range("J14").Value = 0.06760401050708505
$ws.Range("M14").Value = 7.369448
$ws.Range("N14").Value = 14.738896
$ws.Range("O14").Value = 0.7452608427984224
$ws.Range("P14").Value = 0.661061693471796
$ws.Range("Q14").Value = 138.773058304176
$ws.Range("R14").Value = 555.092233216704
$ws.Range("S14").Value = 0.05633012663208906
$ws.Range("T14").Value = 0.04469042167129873
$ws.Range("G15").Value = 18.830862
$ws.Range("H15").Value = 37.661724
$ws.Range("I15").Value = 0.07558444426057843
$ws.Range("J15").Value = 0.06760401050708505
$ws.Range("O15").Value = 0.01116592909756377
$ws.Range("P15").Value = 0.01485661309677453
$ws.Range("Q15").Value = 2.07917824296
$ws.Range("R15").Value = 12.47506945776
$ws.Range("S15").Value = 0.0008439705454923798
$ws.Range("T15").Value = 0.001004366627894042
$ws.Range("G16").Value = 18.830862
$ws.Range("H16").Value = 37.661724
$ws.Range("I16").Value = 0.07558444426057843
$ws.Range("J16").Value = 0.06760401050708505
$ws.Range("M16").Value = 0.084843
$ws.Range("N16").Value = 0.254529
$ws.Range("O16").Value = 0.008580040959044227
$ws.Range("P16").Value = 0.0114160091622658
$ws.Range("Q16").Value = 1.597666824666
$ws.Range("R16").Value = 9.586000947996
$ws.Range("S16").Value = 0.0006485176276223583
$ws.Range("T16").Value = 0.0007717680033547964
$ws.Range("G17").Value = 18.830862
$ws.Range("H17").Value = 37.661724
$ws.Range("I17").Value = 0.07558444426057843
$ws.Range("J17").Value = 0.06760401050708505
$ws.Range("M17").Value = 2.32371
$ws.Range("N17").Value = 6.97113
$ws.Range("O17").Value = 0.2349931871449696
$ws.Range("P17").Value = 0.3126656842691638
$ws.Range("Q17").Value = 43.75746233802001
$ws.Range("R17").Value = 262.54477402812
$ws.Range("S17").Value = 0.01776182945537463
$ws.Range("T17").Value = 0.02113745420453748
$ws.Range("G18").Value = 1.059785333333333
$ws.Range("H18").Value = 3.179356
$ws.Range("I18").Value = 0.004253829986939094
$ws.Range("J18").Value = 0.005707046667055495
$ws.Range("M18").Value = 7.369448
$ws.Range("N18").Value = 14.738896
$ws.Range("O18").Value = 0.7452608427984224
$ws.Range("P18").Value = 0.661061693471796
$ws.Range("Q18").Value = 7.810032905162667
$ws.Range("R18").Value = 46.86019743097601
$ws.Range("S18").Value = 0.003170212921187431
$ws.Range("T18").Value = 0.003772709934446274
$ws.Range("G19").Value = 1.059785333333333
$ws.Range("H19").Value = 3.179356
$ws.Range("I19").Value = 0.004253829986939094
$ws.Range("J19").Value = 0.005707046667055495
$ws.Range("O19").Value = 0.01116592909756377
$ws.Range("P19").Value = 0.01485661309677453
$ws.Range("Q19").Value = 0.1170144312711111
$ws.Range("R19").Value = 1.05312988144
$ws.Range("S19").Value = 0.00004749796402725254
$ws.Range("T19").Value = 0.00008478738425768008
$ws.Range("G20").Value = 1.059785333333333
$ws.Range("H20").Value = 3.179356
$ws.Range("I20").Value = 0.004253829986939094
$ws.Range("J20").Value = 0.005707046667055495
$ws.Range("M20").Value = 0.084843
$ws.Range("N20").Value = 0.254529
$ws.Range("O20").Value = 0.008580040959044227
$ws.Range("P20").Value = 0.0114160091622658
$ws.Range("Q20").Value = 0.089915367036
$ws.Range("R20").Value = 0.809238303324
$ws.Range("S20").Value = 0.00003649803552074799
$ws.Range("T20").Value = 0.00006515169704058403
$ws.Range("G21").Value = 1.059785333333333
$ws.Range("H21").Value = 3.179356
$ws.Range("I21").Value = 0.004253829986939094
$ws.Range("J21").Value = 0.005707046667055495
$ws.Range("M21").Value = 2.32371
$ws.Range("N21").Value = 6.97113
$ws.Range("O21").Value = 0.2349931871449696
$ws.Range("P21").Value = 0.3126656842691638
$ws.Range("Q21").Value = 2.46263377692
$ws.Range("R21").Value = 22.16370399228
$ws.Range("S21").Value = 0.0009996210662036621
$ws.Range("T21").Value = 0.001784397651310957
$ws.Range("G22").Value = 0.146971
$ws.Range("H22").Value = 0.440913
$ws.Range("I22").Value = 0.0005899210220658763
$ws.Range("J22").Value = 0.0007914530700907476
$ws.Range("M22").Value = 7.369448
$ws.Range("N22").Value = 14.738896
$ws.Range("O22").Value = 0.7452608427984224
$ws.Range("P22").Value = 0.661061693471796
$ws.Range("Q22").Value = 1.083095142008
$ws.Range("R22").Value = 6.498570852048
$ws.Range("S22").Value = 0.0004396450380893217
$ws.Range("T22").Value = 0.0005231993068176416
$ws.Range("G23").Value = 0.146971
$ws.Range("H23").Value = 0.440913
$ws.Range("I23").Value = 0.0005899210220658763
$ws.Range("J23").Value = 0.0007914530700907476
$ws.Range("O23").Value = 0.01116592909756377
$ws.Range("P23").Value = 0.01485661309677453
$ws.Range("Q23").Value = 0.01622755801333333
$ws.Range("R23").Value = 0.14604802212
$ws.Range("S23").Value = 0.000006587016305549928
$ws.Range("T23").Value = 0.00001175831204659261
$ws.Range("G24").Value = 0.146971
$ws.Range("H24").Value = 0.440913
$ws.Range("I24").Value = 0.0005899210220658763
$ws.Range("J24").Value = 0.0007914530700907476
$ws.Range("M24").Value = 0.084843
$ws.Range("N24").Value = 0.254529
$ws.Range("O24").Value = 0.008580040959044227
$ws.Range("P24").Value = 0.0114160091622658
$ws.Range("Q24").Value = 0.012469460553
$ws.Range("R24").Value = 0.112225144977
$ws.Range("S24").Value = 0.000005061546531926452
$ws.Range("T24").Value = 0.000009035235499659371
$ws.Range("G25").Value = 0.146971
$ws.Range("H25").Value = 0.440913
$ws.Range("I25").Value = 0.0005899210220658763
$ws.Range("J25").Value = 0.0007914530700907476
$ws.Range("M25").Value = 2.32371
$ws.Range("N25").Value = 6.97113
$ws.Range("O25").Value = 0.2349931871449696
$ws.Range("P25").Value = 0.3126656842691638
$ws.Range("Q25").Value = 0.34151798241
$ws.Range("R25").Value = 3.07366184169
$ws.Range("S25").Value = 0.0001386274211390783
$ws.Range("T25").Value = 0.000247460215726854
